$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 265 (shifts existing rows 265:289 down to 266:290)
$ws.Rows("265:265").Insert()

# Populate the newly inserted row 265 with the new weekly price record
$ws.Range("A265").Value = 11
$ws.Range("B265").Value = "Vega Monumental Concepción"
$ws.Range("C265").Value = "Bíobío"
$ws.Range("D265").Value = 45132
$ws.Range("E265").Value = 8
$ws.Range("F265").Value = "Fruta"
$ws.Range("G265").Value = 100108
$ws.Range("H265").Value = "Tropicales y subtropicales"
$ws.Range("I265").Value = 100108005
$ws.Range("J265").Value = "Piña"
$ws.Range("K265").Value = "Caramelo"
$ws.Range("L265").Value = "Segunda"
$ws.Range("M265").Value = 150
$ws.Range("N265").Value = 19000
$ws.Range("O265").Value = 19000
$ws.Range("P265").Value = 19000
$ws.Range("Q265").Value = "$/caja 14 unidades"
$ws.Range("R265").Value = "Ecuador"
$ws.Range("S265").Value = 1357
$ws.Range("T265").Value = 14
